$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B3").Value = "Viju"
$ws1.Range("A4").Value = "www"
$ws1.Range("C4").Value = "nds.nitin@gmail.com "
$ws1.Rows.Item(8).Delete()

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B27").Select()

$ws1.Activate()
$ws1.Range("B8").Select()
